$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value could be misparsed as a number by Excel's
# auto-detection (e.g. "0.570" -> 0.57, "504.77" stays but could round-trip
# oddly, "0.0210" -> 2.1E-2). Force these to literal text, then restore the
# default "Normal" style so no stray formatting is left behind.
function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '57.171.62'
$ws.Range('E2').Value = '  +0.90%  '
$ws.Range('D3').Value = '2.396.89'
$ws.Range('E3').Value = '  +1.69%  '
$ws.Range('E4').Value = '  +0.08%  '
Set-TextValue 'D5' '504.77'
$ws.Range('E5').Value = '  -1.58%  '
Set-TextValue 'D6' '132.85'
$ws.Range('E6').Value = '  +4.24%  '
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  +0.06%  '
Set-TextValue 'D8' '0.554'
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').Value = '2.410.27'
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('E10').Value = '  +1.19%  '
Set-TextValue 'D11' '0.150'
$ws.Range('E11').Value = '  -1.43%  '
$ws.Range('E12').Value = '  +1.38%  '
Set-TextValue 'D13' '4.58'
$ws.Range('E13').Value = '  -4.71%  '
$ws.Range('D14').Value = '2.828.67'
$ws.Range('E14').Value = '  +2.04%  '
$ws.Range('D15').Value = '57.067.09'
$ws.Range('E15').Value = '  +1.02%  '
Set-TextValue 'D16' '21.84'
$ws.Range('E16').Value = '  +1.97%  '
$ws.Range('E17').Value = '  +2.57%  '
$ws.Range('D18').Value = '2.401.34'
$ws.Range('E18').Value = '  +3.02%  '
$ws.Range('E19').Value = '  -0.49%  '
Set-TextValue 'D20' '309.70'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('E21').Value = '  -0.37%  '
Set-TextValue 'D22' '6.35'
$ws.Range('E22').Value = '  +4.29%  '
Set-TextValue 'D23' '5.87'
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('E24').Value = '  +0.29%  '
Set-TextValue 'D25' '65.20'
$ws.Range('E25').Value = '  +0.12%  '
Set-TextValue 'D26' '0.996'
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  -0.23%  '
Set-TextValue 'D28' '0.376'
$ws.Range('E28').Value = '  -3.05%  '
Set-TextValue 'D29' '7.41'
$ws.Range('E29').Value = '  +3.08%  '
Set-TextValue 'D30' '172.69'
$ws.Range('E30').Value = '  -1.22%  '
$ws.Range('D31').Value = '0.0₃0722'
$ws.Range('E31').Value = '  +1.21%  '
$ws.Range('E32').Value = '  +0.12%  '
Set-TextValue 'D33' '1.12'
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('E34').Value = '  -3.04%  '
$ws.Range('E35').Value = '  +0.10%  '
Set-TextValue 'D36' '0.997'
$ws.Range('E36').Value = '  +0.11%  '
Set-TextValue 'D37' '17.96'
$ws.Range('E37').Value = '  +1.95%  '
$ws.Range('E38').Value = '  +1.22%  '
$ws.Range('E39').Value = '  +3.28%  '
Set-TextValue 'D40' '36.67'
$ws.Range('E40').Value = '  +3.41%  '
Set-TextValue 'D41' '0.802'
$ws.Range('E41').Value = '  -0.83%  '
$ws.Range('E42').Value = '  +1.34%  '
Set-TextValue 'D43' '133.09'
$ws.Range('E43').Value = '  +10.63%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D44' '3.35'
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D45' '4.94'
$ws.Range('E45').Value = '  +1.23%  '
Set-TextValue 'D46' '0.570'
Set-TextValue 'D47' '252.56'
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('E48').Value = '  +0.59%  '
Set-TextValue 'D49' '0.0488'
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '17.04'
$ws.Range('E50').Value = '  +2.84%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D51' '0.0210'
$ws.Range('E51').Value = '  +1.01%  '
